$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format to preserve number-like strings (e.g. "604.83", "1.00")
# exactly as text, matching the original workbook's inline-string cell type.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.375.02"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").Value = "3.487.61"
$ws.Range("E3").Value = "  -2.00%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "604.83"
$ws.Range("E5").Value = "  -2.03%  "
$ws.Range("D6").Value = "150.83"
$ws.Range("E6").Value = "  -1.64%  "
$ws.Range("D7").Value = "3.487.50"
$ws.Range("E7").Value = "  -1.98%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("E10").Value = "  +1.90%  "
$ws.Range("D11").Value = "7.57"
$ws.Range("E11").Value = "  +6.33%  "
$ws.Range("D12").Value = "0.432"
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("D13").Value = "0.0000217"
$ws.Range("E13").Value = "  -2.58%  "
$ws.Range("D14").Value = "32.05"
$ws.Range("E14").Value = "  -0.97%  "
$ws.Range("D15").Value = "4.080.28"
$ws.Range("E15").Value = "  -1.65%  "
$ws.Range("D16").Value = "3.485.83"
$ws.Range("E16").Value = "  -1.64%  "
$ws.Range("D17").Value = "67.405.85"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "6.50"
$ws.Range("D20").Value = "15.38"
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("D21").Value = "9.89"
$ws.Range("E21").Value = "  +1.31%  "
$ws.Range("D22").Value = "445.66"
$ws.Range("E22").Value = "  -0.77%  "
$ws.Range("D23").Value = "0.627"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "78.60"
$ws.Range("E24").Value = "  +1.19%  "
$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").Value = "3.629.66"
$ws.Range("E25").Value = "  -1.60%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -5.49%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "8.66"
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "9.95"
$ws.Range("E29").Value = "  -4.03%  "
$ws.Range("E30").Value = "  -1.38%  "
$ws.Range("E31").Value = "  +1.93%  "
$ws.Range("E32").Value = "  +2.18%  "
$ws.Range("D34").Value = "25.61"
$ws.Range("E34").Value = "  -1.78%  "
$ws.Range("E35").Value = "  -1.86%  "
$ws.Range("D36").Value = "1.86"
$ws.Range("E36").Value = "  -0.69%  "
$ws.Range("D37").Value = "3.485.95"
$ws.Range("E37").Value = "  -1.56%  "
$ws.Range("E38").Value = "  -1.35%  "
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("E40").Value = "  +4.78%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "177.91"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "30.05"
$ws.Range("E46").Value = "  +4.68%  "
$ws.Range("D47").Value = "46.45"
$ws.Range("E47").Value = "  +2.29%  "
$ws.Range("D48").Value = "1.29"
$ws.Range("E48").Value = "  +1.87%  "
$ws.Range("D49").Value = "2.54"
$ws.Range("E49").Value = "  -5.15%  "
$ws.Range("D50").Value = "7.61"
$ws.Range("E50").Value = "  -0.62%  "
$ws.Range("D51").Value = "0.252"
$ws.Range("E51").Value = "  -0.68%  "
